# Aprimorei o visual HTML e CSS das paginas do App
#
# 1) Two existing rows get their "Respondido" status flipped from NO to SIM.
# 2) Three new task rows are appended at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Mark rows 2 and 3 ("Respondido") as answered ---------------------
$ws.Range("I2").Value = "SIM"
$ws.Range("I3").Value = "SIM"

# --- 2) Append three new rows (32-34) -------------------------------------
# Copy the look (bold / centered / bordered) of the existing index column
# onto the new index cells before filling in the values.
$ws.Range("A31").Copy()
$ws.Range("A32:A34").PasteSpecial(-4122)

# Row 32
$ws.Range("A32").Value = 30
$ws.Range("B32").Value = "'12-04-2023"
$ws.Range("B32").ClearFormats()
$ws.Range("C32").Value = "417823 - PREMIUM SAÚDE S.A"
$ws.Range("D32").Value = "12/04/2023  14:44:59"
$ws.Range("E32").Value = 12167081
$ws.Range("F32").Value = 8604741
$ws.Range("G32").Value = "JENIFE BIANCA AMORIM PEREIRA"
$ws.Range("H32").Value = "10 dias úteis"
$ws.Range("I32").Value = "NO"
$ws.Range("J32").Value = "Assistencial"
$ws.Range("K32").Value = "Responder  Detalhes"

# Row 33
$ws.Range("A33").Value = 31
$ws.Range("B33").Value = "'12-04-2023"
$ws.Range("B33").ClearFormats()
$ws.Range("C33").Value = "417823 - PREMIUM SAÚDE S.A"
$ws.Range("D33").Value = "12/04/2023  15:02:14"
$ws.Range("E33").Value = 12167165
$ws.Range("F33").Value = 8604741
$ws.Range("G33").Value = "JUCIMAR AMORIM PEREIRA"
$ws.Range("H33").Value = "10 dias úteis"
$ws.Range("I33").Value = "NO"
$ws.Range("J33").Value = "Assistencial"
$ws.Range("K33").Value = "Responder  Detalhes"

# Row 34
$ws.Range("A34").Value = 32
$ws.Range("B34").Value = "'12-04-2023"
$ws.Range("B34").ClearFormats()
$ws.Range("C34").Value = "417823 - PREMIUM SAÚDE S.A"
$ws.Range("D34").Value = "12/04/2023  17:36:17"
$ws.Range("E34").Value = 12167700
$ws.Range("F34").Value = 8605615
$ws.Range("G34").Value = "VIRGILIO ISMAR SANTOS GARCIA"
$ws.Range("H34").Value = "10 dias úteis"
$ws.Range("I34").Value = "NO"
$ws.Range("J34").Value = "Assistencial"
$ws.Range("K34").Value = "Responder  Detalhes"
